$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto data.
# Price column values are plain text (e.g. "26.536.15", "1.0000") -- force
# text formatting so Excel does not reinterpret them as numbers and lose
# trailing zeros / thousands-style dots / numeric precision.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.536.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.733.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4868"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2667"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06216"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.732.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07045"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.72"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.603"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6099"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.518.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9993"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("E19").Value = "  +4.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.957.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.538"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.753"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.241"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.413"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "108.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.984"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08052"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.692"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04551"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.0000"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.614"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.012"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6374"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9026"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.035"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.398"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01509"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -10.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.426"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3901"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.950"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("E47").Value = "  -1.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05384"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.95%  "
$ws.Range("E49").Value = "  -0.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.783"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.249"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.59%  "
